$wb = $excel.ActiveWorkbook

# --- 1. Template sheet: change A1 text, remove B1 ---
$template = $wb.Worksheets.Item("Template")
$template.Range("A1").Value = "Number of the unit (can include text, if you want)"
$template.Range("B1").Clear()

# --- 2. Unit 0 sheet: update the unit-name description text ---
$unit0 = $wb.Worksheets.Item("Unit 0 - test unit")
$unit0.Range("B2").Value = "Test unit for QA purposes. "

# --- 3. Add the new "Unit 2" sheet after "Unit 1 - add more here!" ---
$unit1 = $wb.Worksheets.Item("Unit 1 - add more here!")
$ws = $wb.Worksheets.Add($null, $unit1)
$ws.Name = "Unit 2 - sentences of necessity"

$ws.Range("A1").Value = "Number of the unit (can include text, if you want)"
$ws.Range("B1").Value = 2
$ws.Range("A2").Value = 'Unit name'
$ws.Range("B2").Value = 'Sentences of necessisity '
$ws.Range("A3").Value = 'Unit description'
$ws.Range("B3").Value = 'Always requires the pronoun.'
$ws.Range("A4").Value = 'Update date'
$ws.Range("B4").Value = '2022.01.30'
$ws.Range("A5").Value = 'Other info'
$ws.Range("B5").Value = 'Created by, notes, etc.'
$ws.Range("A7").Value = 'ID (unique)'
$ws.Range("B7").Value = 'Puhekieli tai kirjakieli'
$ws.Range("C7").Value = 'Type of thing missing'
$ws.Range("D7").Value = 'Verb (infinitive)'
$ws.Range("E7").Value = 'If not a verb, translation or hint here (optional)'
$ws.Range("F7").Value = 'Answer'
$ws.Range("G7").Value = 'Lause (replace verb or phrase with ###)'
$ws.Range("H7").Value = 'English translation of sentence'
$ws.Range("A8").Value = 'U2-1'
$ws.Range("B8").Value = 'kirjakieli'
$ws.Range("C8").Value = 'verb'
$ws.Range("D8").Value = 'haluta'
$ws.Range("E8").Value = 'I - include the proper pronoun!'
$ws.Range("F8").Value = 'Minä haluan'
$ws.Range("G8").Value = '### perunoita!'
$ws.Range("H8").Value = 'I want potatoes!'
$ws.Range("A9").Value = 'U2-2'
$ws.Range("B9").Value = 'kirjakieli'
$ws.Range("C9").Value = 'verb'
$ws.Range("D9").Value = 'tarvita'
$ws.Range("E9").Value = 'I - include the proper pronoun!'
$ws.Range("F9").Value = 'Minä tarvitsen'
$ws.Range("G9").Value = '### perunoita!'
$ws.Range("H9").Value = 'I need potatoes!'
$ws.Range("A10").Value = 'U2-3'
$ws.Range("B10").Value = 'kirjakieli'
$ws.Range("C10").Value = 'verb'
$ws.Range("D10").Value = 'tarvita'
$ws.Range("E10").Value = 'You (singular) - include the proper pronoun'
$ws.Range("F10").Value = 'Sinä tarvitset'
$ws.Range("G10").Value = '### perunoita!'
$ws.Range("H10").Value = 'You (singular) need potatoes'
$ws.Range("A11").Value = 'U2-4'
$ws.Range("B11").Value = 'kirjakieli'
$ws.Range("C11").Value = 'verb'
$ws.Range("D11").Value = 'täytyä'
$ws.Range("E11").Value = 'Michelle - include pronoun'
$ws.Range("F11").Value = 'Michellen täytyy'
$ws.Range("G11").Value = '### saada perunat.'
$ws.Range("H11").Value = 'Michelle needs to get potatoes'
$ws.Range("A12").Value = 'U2-5'
$ws.Range("B12").Value = 'kirjakieli'
$ws.Range("C12").Value = 'verb'
$ws.Range("D12").Value = 'täytyä'
$ws.Range("E12").Value = 'Kissa - include pronoun'
$ws.Range("F12").Value = 'Kissan ei tarvitse'
$ws.Range("G12").Value = '### syödä perunaa.'
$ws.Range("H12").Value = 'The cat does not need to eat a potato.'
$ws.Range("A13").Value = 'U2-6'
$ws.Range("B13").Value = 'kirjakieli'
$ws.Range("C13").Value = 'verb'
$ws.Range("D13").Value = 'pitää'
$ws.Range("E13").Value = 'We - include the pronoun'
$ws.Range("F13").Value = 'Meidän pitää'
$ws.Range("G13").Value = '### kokata perunoita'
$ws.Range("H13").Value = 'We need to cook potatoes'
$ws.Range("A14").Value = 'U2-7'
$ws.Range("B14").Value = 'kirjakieli'
$ws.Range("C14").Value = 'verb'
$ws.Range("D14").Value = 'pitää'
$ws.Range("E14").Value = 'Kaupa - include pronoun.  Note the "should"!'
$ws.Range("F14").Value = 'Kaupan pitäisi'
$ws.Range("G14").Value = '### myydä perunoita'
$ws.Range("H14").Value = 'The store should sell potatoes'
$ws.Range("A15").Value = 'U2-8'
$ws.Range("B15").Value = 'kirjakieli'
$ws.Range("C15").Value = 'verb'
$ws.Range("D15").Value = 'pitää'
$ws.Range("E15").Value = 'You (plural) - include pronoun'
$ws.Range("F15").Value = 'Teidän ei pitäisi'
$ws.Range("G15").Value = '### unohtaa perunoita'
$ws.Range("H15").Value = 'You (plural) should not forget potatoes'
$ws.Range("A16").Value = 'U2-9'
$ws.Range("B16").Value = 'kirjakieli'
$ws.Range("C16").Value = 'verb'
$ws.Range("D16").Value = 'olla'
$ws.Range("E16").Value = 'They - include pronoun (and maybe you need another word too?)'
$ws.Range("F16").Value = 'Heidän on pakko'
$ws.Range("G16").Value = '### ostaa perunoita'
$ws.Range("H16").Value = 'They must buy potaotes (it is the law!)'
$ws.Range("A17").Value = 'U2-10'
$ws.Range("B17").Value = 'kirjakieli'
$ws.Range("C17").Value = 'verb'
$ws.Range("D17").Value = 'olla'
$ws.Range("E17").Value = 'We - include pronoun (and maybe another word?)'
$ws.Range("F17").Value = 'Meidän olisi hyvä'
$ws.Range("G17").Value = '### pyöräillä järvin ympäri'
$ws.Range("H17").Value = 'It would be a good idea for us to ride our bikes around the lake.'
$ws.Range("A18").Value = 'U2-11'
$ws.Range("B18").Value = 'kirjakieli'
$ws.Range("C18").Value = 'verb'
$ws.Range("D18").Value = 'olla'
$ws.Range("E18").Value = 'We - include pronoun (and maybe another word?)'
$ws.Range("F18").Value = 'Meidän on pakko'
$ws.Range("G18").Value = '### pyöräillä pyöriäme'
$ws.Range("H18").Value = 'We must ride bicycles'
$ws.Range("A19").Value = 'U2-12'
$ws.Range("B19").Value = 'kirjakieli'
$ws.Range("C19").Value = 'verb'
$ws.Range("D19").Value = 'olla'
$ws.Range("E19").Value = 'I - include pronoun (and maybe an extra word?)'
$ws.Range("F19").Value = 'Minun olisi hyvä'
$ws.Range("G19").Value = '### kävellä joka päivä'
$ws.Range("H19").Value = 'It would be a good idea for me to walk every day'
$ws.Range("A20").Value = 'U2-13'
$ws.Range("B20").Value = 'kirjakieli'
$ws.Range("C20").Value = 'verb'
$ws.Range("D20").Value = 'olla'
$ws.Range("E20").Value = 'You (singular) - include pronoun (and maybe another word?)'
$ws.Range("F20").Value = 'Sinun ei ole pakko'
$ws.Range("G20").Value = '### laulaa'
$ws.Range("H20").Value = 'You (singular) must not sing'
$ws.Range("A21").Value = 'U2-14'
$ws.Range("B21").Value = 'kirjakieli'
$ws.Range("C21").Value = 'verb'
$ws.Range("D21").Value = 'olla'
$ws.Range("E21").Value = 'Matti - include pronoun (and maybe another word?)'
$ws.Range("F21").Value = 'Mattin ei ole pakko'
$ws.Range("G21").Value = '### juosta koska hänellä on polvi kipeä'
$ws.Range("H21").Value = 'Matti must not run because he has a hurt knee'
$ws.Range("A22").Value = 'U2-15'
$ws.Range("B22").Value = 'kirjakieli'
$ws.Range("C22").Value = 'verb'
$ws.Range("D22").Value = 'olla'
$ws.Range("E22").Value = 'Pekka - include pronoun'
$ws.Range("F22").Value = 'Pekkan ei ole pakko'
$ws.Range("G22").Value = '### ajaa autoa koska hän on pyöräilijäksi'
$ws.Range("H22").Value = 'Pekka must not drive because he is the cyclist of the year'
$ws.Range("A23").Value = 'U2-16'
$ws.Range("B23").Value = 'kirjakieli'
$ws.Range("C23").Value = 'verb'
$ws.Range("D23").Value = 'olla'
$ws.Range("E23").Value = 'She - inclue pronoun (and maybe an extra word?)'
$ws.Range("F23").Value = 'Hänen olisi hyvä'
$ws.Range("G23").Value = '### juoda teeta koska hänellä on kylmä'
$ws.Range("H23").Value = 'It would be a good idea for her to drink tea, becuase she is cold.'
$ws.Range("A24").Value = 'U2-17'
$ws.Range("B24").Value = 'kirjakieli'
$ws.Range("C24").Value = 'verb'
$ws.Range("D24").Value = 'täytyä'
$ws.Range("E24").Value = 'Oppetaja - include pronoun'
$ws.Range("F24").Value = 'Oppetajan täytyy'
$ws.Range("G24").Value = '### nukua'
$ws.Range("H24").Value = 'The student needs to sleep'
$ws.Range("A25").Value = 'U2-18'
$ws.Range("B25").Value = 'kirjakieli'
$ws.Range("C25").Value = 'verb'
$ws.Range("D25").Value = 'täytyä'
$ws.Range("E25").Value = 'Lucca and Kira - include pronoun'
$ws.Range("F25").Value = 'Luccan ja Kiran täytyy'
$ws.Range("G25").Value = '### syödä lihaa.'
$ws.Range("H25").Value = 'Lucca and Kira need to eat meat'
$ws.Range("A26").Value = 'U2-19'
$ws.Range("B26").Value = 'kirjakieli'
$ws.Range("C26").Value = 'verb'
$ws.Range("D26").Value = 'pitää'
$ws.Range("E26").Value = 'They - include pronoun  '
$ws.Range("F26").Value = 'Heidän pitää'
$ws.Range("G26").Value = '### ostaa kissanruoka'
$ws.Range("H26").Value = 'They need to buy cat food'
$ws.Range("A27").Value = 'U2-20'
$ws.Range("B27").Value = 'kirjakieli'
$ws.Range("C27").Value = 'verb'
$ws.Range("D27").Value = 'pitää'
$ws.Range("E27").Value = 'Bertie - include pronoun'
$ws.Range("F27").Value = 'Bertien ei pitäisi'
$ws.Range("G27").Value = '### puraista Lucca'
$ws.Range("H27").Value = 'Bertie should not bite Luuca'
$ws.Range("A28").Value = 'U2-21'
$ws.Range("B28").Value = 'kirjakieli'
$ws.Range("C28").Value = 'verb'
$ws.Range("D28").Value = 'pitää'
$ws.Range("E28").Value = 'We - include pronoun'
$ws.Range("F28").Value = 'Meidän pitäisi'
$ws.Range("G28").Value = '### siivota'
$ws.Range("H28").Value = 'We should clean'
$ws.Range("A29").Value = 'U2-22'
$ws.Range("B29").Value = 'kirjakieli'
$ws.Range("C29").Value = 'verb'
$ws.Range("D29").Value = 'tarvita'
$ws.Range("E29").Value = 'We - include pronoun'
$ws.Range("F29").Value = 'Me tarvitsimme'
$ws.Range("G29").Value = '### perunoita!'
$ws.Range("H29").Value = 'We need potatoes!'
$ws.Range("A30").Value = 'U2-23'
$ws.Range("B30").Value = 'kirjakieli'
$ws.Range("C30").Value = 'verb'
$ws.Range("D30").Value = 'pitää'
$ws.Range("E30").Value = 'We - include pronoun'
$ws.Range("F30").Value = 'Meidän ei tarvitse'
$ws.Range("G30").Value = '### ajaa autoa suomessa.'
$ws.Range("H30").Value = 'We do not need to drive in Finland.'
$ws.Range("A31").Value = 'U2-24'
$ws.Range("B31").Value = 'kirjakieli'
$ws.Range("C31").Value = 'verb'
$ws.Range("D31").Value = 'pitää'
$ws.Range("E31").Value = 'Kaikki - include pronoun'
$ws.Range("F31").Value = 'Kaikkien pitäisi'
$ws.Range("G31").Value = '### puhua suomea'
$ws.Range("H31").Value = 'Everyone should speak Finnish.'
$ws.Range("A32").Value = 'U2-25'
$ws.Range("B32").Value = 'kirjakieli'
$ws.Range("C32").Value = 'verb'
$ws.Range("D32").Value = 'pitää'
$ws.Range("E32").Value = 'Michelle - include pronoun'
$ws.Range("F32").Value = 'Michellen pitää'
$ws.Range("G32").Value = '### lisätä maitoa ja sokeria kahviinsa'
$ws.Range("H32").Value = 'Michelle needs to add milk and sugar to her coffee'
$ws.Range("A33").Value = 'U2-26'
$ws.Range("B33").Value = 'kirjakieli'
$ws.Range("C33").Value = 'verb'
$ws.Range("D33").Value = 'täytyä'
$ws.Range("E33").Value = 'She - include pronoun'
$ws.Range("F33").Value = 'Hänen täytyy'
$ws.Range("G33").Value = '### korjata pyoränsä'
$ws.Range("H33").Value = 'She needs to fix her bike'
$ws.Range("A34").Value = 'U2-27'
$ws.Range("B34").Value = 'kirjakieli'
$ws.Range("C34").Value = 'verb'
$ws.Range("D34").Value = 'täytyä'
$ws.Range("E34").Value = 'They - include pronoun'
$ws.Range("F34").Value = 'Heidän täytyy'
$ws.Range("G34").Value = '### lainata kirja kirjastosta'
$ws.Range("H34").Value = 'They need to borrow a book from the library'
$ws.Range("A35").Value = 'U2-28'
$ws.Range("B35").Value = 'kirjakieli'
$ws.Range("C35").Value = 'verb'
$ws.Range("D35").Value = 'tarvita'
$ws.Range("E35").Value = 'They - include pronoune'
$ws.Range("F35").Value = 'He tarvitsevat'
$ws.Range("G35").Value = '### kirja kirjastosta'
$ws.Range("H35").Value = 'They need a book from the library'
$ws.Range("A36").Value = 'U2-29'
$ws.Range("B36").Value = 'kirjakieli'
$ws.Range("C36").Value = 'verb'
$ws.Range("D36").Value = 'olla'
$ws.Range("E36").Value = 'Kaikki - include pronoun (and another word?)'
$ws.Range("F36").Value = 'Kaikkien on pakko'
$ws.Range("G36").Value = '### olla hiljaa yöllä rakennuksessa '
$ws.Range("H36").Value = 'Everyone must be quiet at night in the building'
$ws.Range("A37").Value = 'U2-30'
$ws.Range("B37").Value = 'kirjakieli'
$ws.Range("C37").Value = 'verb'
$ws.Range("D37").Value = 'pitää'
$ws.Range("E37").Value = 'You (plural) - include pronoun'
$ws.Range("F37").Value = 'Teidän pitää'
$ws.Range("G37").Value = '### joskus katsoa ulos'
$ws.Range("H37").Value = 'You (plural) need to look outside sometime'

# Column widths to match the other unit sheets
$ws.Columns.Item(1).ColumnWidth = 13.86
$ws.Columns.Item(2).ColumnWidth = 14.0
$ws.Columns.Item(5).ColumnWidth = 24.29
$ws.Columns.Item(6).ColumnWidth = 24.29
$ws.Columns.Item(7).ColumnWidth = 42.0
$ws.Columns.Item(8).ColumnWidth = 46.43
